$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new data row at 22 (pushes old row22..28 down to 23..29) ---
$ws.Rows("22:22").Insert()

# Copy formatting (borders, fill, number format, font) from row 21 (a normal
# data row) into the freshly inserted row 22 so it matches its neighbours.
$ws.Range("B21:J21").Copy()
$ws.Range("B22:J22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Header block updates ---
$ws.Range("E11").Value2 = 387192
$ws.Range("C13").Value2 = 8

# Header row (row 15) relabels: "Novedad de Ingreso"/"Novedad de Retiro" swap
# with "Observaciones" moving between them.
$ws.Range("H15").Value2 = "Novedad de Retiro"
$ws.Range("I15").Value2 = "Observaciones"
$ws.Range("J15").Value2 = "Novedad de Ingreso"

# --- Worker rows 16-23 (new data set) ---
# Row 16: EVELISBETH moves to the top of the list, same 2505 period.
$ws.Range("B16").Value2 = "CC"
$ws.Range("C16").Value2 = "1047522475"
$ws.Range("D16").Value2 = "EVELISBETH DEL CARMEN LOBO PINO"
$ws.Range("E16").Value2 = "2505"
$ws.Range("F16").Value2 = 17082
$ws.Range("G16").Value2 = 1423500

# Row 17: YOHANA AVILA ROBLES, period moves 2507 -> 2508.
$ws.Range("B17").Value2 = "CC"
$ws.Range("C17").Value2 = "1143333388"
$ws.Range("D17").Value2 = "YOHANA AVILA ROBLES"
$ws.Range("E17").Value2 = "2508"
$ws.Range("F17").Value2 = 56940
$ws.Range("G17").Value2 = 1423500

# Row 18: DAYANA NEWBALL TINOCO, period -> 2508.
$ws.Range("B18").Value2 = "CC"
$ws.Range("C18").Value2 = "1047426266"
$ws.Range("D18").Value2 = "DAYANA NEWBALL TINOCO"
$ws.Range("E18").Value2 = "2508"
$ws.Range("F18").Value2 = 56940
$ws.Range("G18").Value2 = 1423500

# Row 19: HAROLD FRANCISCO ROMERO BLANQUICET, period -> 2508.
$ws.Range("B19").Value2 = "CC"
$ws.Range("C19").Value2 = "1007976089"
$ws.Range("D19").Value2 = "HAROLD FRANCISCO ROMERO BLANQUICET"
$ws.Range("E19").Value2 = "2508"
$ws.Range("F19").Value2 = 56940
$ws.Range("G19").Value2 = 1423500

# Row 20: YERSON JIMENEZ GUZMAN, period -> 2508.
$ws.Range("B20").Value2 = "CC"
$ws.Range("C20").Value2 = "1049830970"
$ws.Range("D20").Value2 = "YERSON JIMENEZ GUZMAN"
$ws.Range("E20").Value2 = "2508"
$ws.Range("F20").Value2 = 56940
$ws.Range("G20").Value2 = 1423500

# Row 21: JOSE DEL CARMEN VILLAMIZAR CHACON, period -> 2508.
$ws.Range("B21").Value2 = "CC"
$ws.Range("C21").Value2 = "1050978914"
$ws.Range("D21").Value2 = "JOSE DEL CARMEN VILLAMIZAR CHACON"
$ws.Range("E21").Value2 = "2508"
$ws.Range("F21").Value2 = 56940
$ws.Range("G21").Value2 = 1423500

# Row 22 (new): MARIA ERMELINDA GUTIERREZ ROJAS, newly added worker.
$ws.Range("B22").Value2 = "CC"
$ws.Range("C22").Value2 = "1130659488"
$ws.Range("D22").Value2 = "MARIA ERMELINDA GUTIERREZ ROJAS"
$ws.Range("E22").Value2 = "2508"
$ws.Range("F22").Value2 = 28470
$ws.Range("G22").Value2 = 1423500

# Row 23 (was row 22): KATERIN JUDITH ESCORCIA PAJARO, period -> 2508.
$ws.Range("B23").Value2 = "CC"
$ws.Range("C23").Value2 = "1007210618"
$ws.Range("D23").Value2 = "KATERIN JUDITH ESCORCIA PAJARO"
$ws.Range("E23").Value2 = "2508"
$ws.Range("F23").Value2 = 56940
$ws.Range("G23").Value2 = 1423500
